# "Daily 100 Error Counts" upload refresh:
# two more days of data were appended (10/29 and 10/28), and the sheet was
# left scrolled down with the newest row selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 16: 2025-10-29 (date serial 45959)
$ws.Range("A16").Value = 45959
$ws.Range("B16").Value = 641
$ws.Range("C16").Value = 16
$ws.Range("D16").Value = 625

# New row 17: 2025-10-28 (date serial 45958)
$ws.Range("A17").Value = 45958
$ws.Range("B17").Value = 648
$ws.Range("C17").Value = 23
$ws.Range("D17").Value = 625

# Leave the sheet scrolled so row 8 is at the top, with the last entered
# row selected, matching where the user's cursor ended up after typing.
$ws.Range("A17:D17").Select()
$excel.ActiveWindow.ScrollRow = 8
